$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right after the H1 title at the top of the document.
# ------------------------------------------------------------------
$metaRange = $d.Content
$metaRange.Find.ClearFormatting()
$found = $metaRange.Find.Execute("Meta description: Read our review of Cashpot Kegs, a unique slot game with exciting bonus features and a Cashpot jackpot. Play for free and discover your chance to win big.", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
if ($found) {
    # Also swallow the trailing paragraph mark so no blank paragraph is left behind.
    $delRange = $d.Range($metaRange.Start, $metaRange.End + 1)
    $delRange.Delete()
}

# ------------------------------------------------------------------
# 2) Insert a brand-new bold paragraph ("Play Cashpot Kegs Free: Unique
#    Slot Game with Exciting Features") right before the final
#    paragraph. We build it via raw OOXML (Range.InsertXML) so the run
#    layout matches exactly: a leading empty run followed by a single
#    bold run, with no stray direct-formatting overrides.
# ------------------------------------------------------------------
$boldText = "Play Cashpot Kegs Free: Unique Slot Game with Exciting Features"
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertStart = $lastPara.Range.Start
$insertPoint = $d.Range($insertStart, $insertStart)

$openXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>' + $boldText + '</w:t></w:r></w:p></w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($openXml)

# The inserted run lands merged onto the front of the old final
# paragraph; split it back into its own paragraph right after the
# bold text.
$splitPos = $insertStart + $boldText.Length
$splitPoint = $d.Range($splitPos, $splitPos)
$splitPoint.InsertParagraphBefore()

# ------------------------------------------------------------------
# 3) Replace the (now last) paragraph's italic text — the old AI-image
#    prompt — with the former meta-description body, keeping its
#    italic formatting untouched.
# ------------------------------------------------------------------
$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute("Create an appealing feature image for Cashpot Kegs Please create an image that represents the Cashpot Kegs game in a cartoon style. The image should feature a happy Maya warrior wearing glasses to make it more interesting. The Maya warrior should have a smile on his face and be standing in front of a stack of gold coins with a jug of beer beside him. In the background, there should be an ancient temple with the sun setting behind it. The image should be colorful and vibrant to attract players to the game. Please make sure the image is high-resolution so that it can be used across various platforms to promote the game.", $true, $false, $false, $false, $false, $true, 1, $false, "Read our review of Cashpot Kegs, a unique slot game with exciting bonus features and a Cashpot jackpot. Play for free and discover your chance to win big.", 2)
